$d = $word.ActiveDocument

# --- Change 1: the empty paragraph right before "Având cluster-ul gata..."
# gets an explicit Romanian language tag on its (empty) run.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $trimmed = $p.Range.Text.Trim()
    if ($trimmed -eq "" -and $i -lt $d.Paragraphs.Count) {
        $nextText = $d.Paragraphs($i + 1).Range.Text
        if ($nextText -like "*cluster-ul gata*") {
            $targetPara = $p
            break
        }
    }
}
if ($targetPara -ne $null) {
    $targetPara.Range.LanguageID = "ro-RO"
}

# --- Change 2: fix the diacritics in the final paragraph's text
# ("usoara" -> "ușoară").
$d.Content.Find.Execute("usoara", $true, $false, $false, $false, $false,
                         $true, 1, $false, "ușoară", 2)

# --- Change 3: append two new, empty Normal paragraphs (lang ro-RO on the
# paragraph mark, plain/empty run) after the last paragraph in the body.
$newEmptyParaXml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:rPr><w:lang w:val="ro-RO"/></w:rPr></w:pPr><w:r><w:rPr/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$insertionPoint = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$insertionPoint.InsertXML($newEmptyParaXml)

# Word's "smart" paste drops explicit spacing that already equals the
# "Normal" style default (before=0/after=160) - so re-assert it on both
# freshly inserted paragraphs to force it back into the OOXML explicitly,
# matching the rest of the document's paragraphs.
$count = $d.Paragraphs.Count
$d.Paragraphs($count).SpaceBefore = 0
$d.Paragraphs($count).SpaceAfter = 8
$d.Paragraphs($count - 1).SpaceBefore = 0
$d.Paragraphs($count - 1).SpaceAfter = 8
